$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Components")

# The Hope RF RFM98W 433 LoRa module is now sourced from HABSupplies
# instead of AnArduino/eBay - update price, supplier, notes and URL.
$ws.Range("C5").Value = 5.99
$ws.Range("E5").Value = "HABSupplies"
$ws.Range("G5").Value = "LoRa 433mhz module. Gnd and CS module"
$ws.Range("H5").Value = "http://ava.upuaut.net/store/index.php?route=product/product&product_id=110"

# Leave the cursor/selection on the cell that was last edited.
$ws.Range("H5").Select()
